$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The "Meta description" paragraph (currently paragraph 2, right after
# the H1 title) needs to move down to become the description line that
# sits just above the final "Prompt: ..." paragraph, and that trailing
# paragraph's text needs to become the description text (keeping its
# italic formatting). The old bold "Meta description" paragraph also
# gets re-purposed into a new bold heading-style line.
# ------------------------------------------------------------------

# Step 1: grab paragraph 2 ("Meta description: ...") as formatted text,
# including its own paragraph mark, so its run/formatting structure is
# preserved verbatim.
$metaSource = $d.Paragraphs.Item(2)
$metaFullRange = $d.Range($metaSource.Range.Start, $metaSource.Range.End)
$metaFormatted = $metaFullRange.FormattedText

# Step 2: insert a copy of that paragraph right before the very last
# paragraph in the document (the one holding the old image "Prompt").
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertionPoint.FormattedText = $metaFormatted

# Step 3: remove the original "Meta description" paragraph from the top
# of the document.
$metaOriginal = $d.Paragraphs.Item(2)
$metaOriginal.Range.Delete()

# Step 4: turn the relocated copy into the new bold heading line
# ("Play Egyptian Riches for Free - Exciting Ancient Egypt-themed Slot").
$newHeadingText = "Play Egyptian Riches for Free - Exciting Ancient Egypt-themed Slot"
$relocatedParaIndex = $d.Paragraphs.Count - 1
$relocatedPara = $d.Paragraphs.Item($relocatedParaIndex)
$oldMetaText = "Meta description: Discover the ancient world of Egypt with Egyptian Riches slot game. Play for free and enjoy themed symbols, bonus features, and 20 paylines."
$relocatedPara.Range.Find.Execute($oldMetaText, $true, $false, $false, $false, $false, $true, 1, $false, $newHeadingText, 2)
$relocatedTextRange = $d.Range($relocatedPara.Range.Start, $relocatedPara.Range.Start + $newHeadingText.Length)
$relocatedTextRange.Font.Bold = 1

# Step 5: replace the old image-prompt text in the final paragraph with
# the meta description text, keeping the paragraph's italic formatting.
$oldPromptText = "Prompt: Create a feature image for Egyptian Riches that captures the excitement and theme of the game. The image should be in a cartoon style and feature a happy Maya warrior wearing glasses. The warrior should be holding up a bag of treasure with hieroglyphics in the background. The image should be colorful and eye-catching, with a playful and adventurous tone. The overall design should entice potential players to join in on the fun and excitement of this Ancient Egyptian themed slot game."
$newDescriptionText = "Discover the ancient world of Egypt with Egyptian Riches slot game. Play for free and enjoy themed symbols, bonus features, and 20 paylines."
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.Range.Find.Execute($oldPromptText, $true, $false, $false, $false, $false, $true, 1, $false, $newDescriptionText, 2)
